# Arabic translation pass for "Email 5-1 [TEMPLATE] Partner email – invite revoked"
$d = $word.ActiveDocument

# "English" appears twice in the body (link label + standalone language heading);
# both map to the same Arabic translation, so a single ReplaceAll covers them.
$d.Content.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "الإنجليزية", 2)

$d.Content.Find.Execute(" / Portuguese / French / Thai / Vietnamese / Spanish", $true, $false, $false, $false, $false, $true, 1, $false, " /البرتغالية/الفرنسية/التايلندية/الفيتنامية/الإسبانية", 2)

$d.Content.Find.Execute("Brief", $true, $false, $false, $false, $false, $true, 1, $false, "المضمون", 2)

$d.Content.Find.Execute("An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. It will be sent via customer.io", $true, $false, $false, $false, $false, $true, 1, $false, "An email sent to partners in the target country who RSVPed yes but didn’t submit their documents by the deadline. We will be revoking their invites. سيتم إرسالها عبر customer.io", 2)

$d.Content.Find.Execute("Target audience", $true, $false, $false, $false, $false, $true, 1, $false, "الجمهور المستهدف", 2)

$d.Content.Find.Execute("We didn’t receive your documents on time", $true, $false, $false, $false, $false, $true, 1, $false, "لم نستلم مستنداتك في الموعد المحدد", 2)

$d.Content.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, "مرحبًا  ", 2)

# ", " is not unique document-wide (it also occurs mid-sentence elsewhere), so scope the
# search to the "Hi [PARTNER NAME], " greeting paragraph (identified by the placeholder,
# which this pass doesn't touch) where the lone trailing ", " run lives.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pGreeting = $d.Paragraphs($i)
    if ($pGreeting.Range.Text.Contains("[PARTNER NAME]")) {
        $pGreeting.Range.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, ",، ", 2)
        break
    }
}

$d.Content.Find.Execute("If you have any questions, please contact us via ", $true, $false, $false, $false, $false, $true, 1, $false, "إذا كانت لديك أي أسئلة، فاتصل بنا:  ", 2)

$d.Content.Find.Execute("live chat", $true, $false, $false, $false, $false, $true, 1, $false, "الدردشة الحية", 2)

$d.Content.Find.Execute("If you have any questions, please contact your country manager, ", $true, $false, $false, $false, $false, $true, 1, $false, "إذا كانت لديك أي أسئلة، فيُرجى الاتصال بمدير بلدك  ", 2)

# Comments live in a separate story; update each comment's range directly.
foreach ($c in $d.Comments) {
    if ($c.Range.Text -eq "choose either one") {
        $c.Range.Text = "اختر أيًا منهما"
    }
}
